$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 2496
$ws1.Range("F3").Value = 533
$ws1.Range("F4").Value = 446
$ws1.Range("F5").Value = 274
$ws1.Range("F6").Value = 167
$ws1.Range("F7").Value = 431
$ws1.Range("F8").Value = 1160
$ws1.Range("F9").Value = 521
$ws1.Range("F10").Value = 275
$ws1.Range("F11").Value = 106
$ws1.Range("F12").Value = 329
$ws1.Range("F13").Value = 5357
$ws1.Range("F15").Value = 1564
$ws1.Range("F16").Value = 3923
$ws1.Range("F17").Value = 379
$ws1.Range("F18").Value = 237
$ws1.Range("F19").Value = 304
$ws1.Range("F20").Value = 4386
$ws1.Range("F21").Value = 5792
$ws1.Range("F24").Value = 627
$ws1.Range("F25").Value = 3604
$ws1.Range("F26").Value = 448
$ws1.Range("F28").Value = 171
$ws1.Range("F29").Value = 112
$ws1.Range("F30").Value = 947
$ws1.Range("F31").Value = 1313
$ws1.Range("F32").Value = 117
$ws1.Range("F33").Value = 158
$ws1.Range("F34").Value = 1533
$ws1.Range("F35").Value = 177
$ws1.Range("F36").Value = 1588
$ws1.Range("F37").Value = 133
$ws1.Range("F38").Value = 1034
$ws1.Range("F39").Value = 27
$ws1.Range("F41").Value = 586
$ws1.Range("F43").Value = 165
$ws1.Range("F44").Value = 2689
$ws1.Range("F45").Value = 114
$ws1.Range("F46").Value = 237
$ws1.Range("F47").Value = 392
$ws1.Range("F49").Value = 3833

# Row 40 also flips G40 from a numeric min-price to "已售罄" (sold out)
$ws1.Range("F40").Value = 1343
$ws1.Range("G40").Value = "已售罄"

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 1141
$ws2.Range("F22").Value = 62

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 3571

# Sheet 4: 全部类型 (All types) -- union of the above, independently maintained
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 3571
$ws4.Range("F3").Value = 2496
$ws4.Range("F4").Value = 533
$ws4.Range("F5").Value = 446
$ws4.Range("F6").Value = 274
$ws4.Range("F7").Value = 1141
$ws4.Range("F8").Value = 167
$ws4.Range("F9").Value = 431
$ws4.Range("F10").Value = 1160
$ws4.Range("F11").Value = 521
$ws4.Range("F12").Value = 275
$ws4.Range("F13").Value = 106
$ws4.Range("F14").Value = 329
$ws4.Range("F15").Value = 5358
$ws4.Range("F17").Value = 1565
$ws4.Range("F18").Value = 4387
$ws4.Range("F19").Value = 5792
$ws4.Range("F22").Value = 627
$ws4.Range("F23").Value = 3604
$ws4.Range("F24").Value = 448
$ws4.Range("F26").Value = 171
$ws4.Range("F27").Value = 112
$ws4.Range("F28").Value = 947
$ws4.Range("F29").Value = 1313
$ws4.Range("F30").Value = 117
$ws4.Range("F31").Value = 158
$ws4.Range("F32").Value = 1533
$ws4.Range("F33").Value = 177
$ws4.Range("F34").Value = 1588
$ws4.Range("F36").Value = 1034
$ws4.Range("F38").Value = 586
$ws4.Range("F42").Value = 62
$ws4.Range("F43").Value = 2689
$ws4.Range("F45").Value = 114
$ws4.Range("F46").Value = 237
$ws4.Range("F47").Value = 392
$ws4.Range("F49").Value = 3833

